# "updated End Date against each Bug"
# Populate column G (End Date) for the bug rows that were still missing it,
# and widen the assignee on bug #6 (row 8) to include Mukesh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Use the existing Start Date cell in F3 (date-formatted, m/d/yyyy) as the
# format template for the new End Date cells in column G, so the new cells
# reuse the workbook's existing date style instead of creating a new one.
$fmtSource = $ws.Range("F3")
$fmtSource.Copy()

$endDates = @{
    3  = 42645
    5  = 42645
    6  = 42676
    7  = 42615
    8  = 42645
    9  = 42645
    10 = 42676
    11 = 42706
    13 = 42645
    14 = 42645
}

foreach ($row in 3, 5, 6, 7, 8, 9, 10, 11, 13, 14) {
    $cell = $ws.Range("G$row")
    $cell.PasteSpecial(-4122)   # xlPasteFormats
    $cell.Value = $endDates[$row]
}

$excel.CutCopyMode = 0

# Bug #6 (row 8) is now also assigned to Mukesh
$ws.Range("C8").Value = "Shruti/Mukesh"

# Update the saved view: scrolled down so row 10 is at the top, G8 selected
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("G8").Select()
